# Update the "Eventos" sheet with the latest daily events export.
# Rows 2-6 are refreshed with new event data and the former row 7 is
# removed entirely (dataset shrinks from A1:L7 to A1:L6).
#
# Columns A-G, K, L hold text values (even when they look numeric, e.g.
# driverId/vehicleId/Unidad). A leading apostrophe forces Excel to store
# them as text instead of auto-converting to numbers, and resetting the
# cell Style back to "Normal" afterwards avoids leaving a stray
# text/quote-prefix number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

# --- Row 2 -----------------------------------------------------------
Set-TextCell "A2" "281474991265569-1738880263843"
Set-TextCell "B2" "Harsh Brake"
Set-TextCell "C2" "2025-02-06T16:17:43.843"
Set-TextCell "D2" "281474991265569"
Set-TextCell "E2" "103"
Set-TextCell "F2" "52215661"
Set-TextCell "G2" "CARLOS ALBERTO JIMENEZ"
$ws.Range("H2").Value = 20.60917552
$ws.Range("I2").Value = -103.42308905
$ws.Range("J2").Value = 0.758480966091156
Set-TextCell "K2" "No video URL"
Set-TextCell "L2" "No video URL"

# --- Row 3 -----------------------------------------------------------
Set-TextCell "A3" "281474991205962-1738874308509"
Set-TextCell "B3" "Harsh Brake"
Set-TextCell "C3" "2025-02-06T14:38:28.509"
Set-TextCell "D3" "281474991205962"
Set-TextCell "E3" "130"
Set-TextCell "F3" "52211469"
Set-TextCell "G3" "DANIEL MERCADO"
$ws.Range("H3").Value = 20.67301142
$ws.Range("I3").Value = -103.33031378
$ws.Range("J3").Value = 0.7648777961730957
Set-TextCell "K3" "No video URL"
Set-TextCell "L3" "No video URL"

# --- Row 4 -----------------------------------------------------------
Set-TextCell "A4" "281474991265569-1738872005633"
Set-TextCell "B4" "Harsh Brake"
Set-TextCell "C4" "2025-02-06T14:00:05.633"
Set-TextCell "D4" "281474991265569"
Set-TextCell "E4" "103"
Set-TextCell "F4" "52215661"
Set-TextCell "G4" "CARLOS ALBERTO JIMENEZ"
$ws.Range("H4").Value = 20.66851241
$ws.Range("I4").Value = -103.37152066
$ws.Range("J4").Value = 0.7957859635353088
Set-TextCell "K4" "No video URL"
Set-TextCell "L4" "No video URL"

# --- Row 5 -----------------------------------------------------------
Set-TextCell "A5" "281474991265569-1738871962133"
Set-TextCell "B5" "Harsh Brake"
Set-TextCell "C5" "2025-02-06T13:59:22.133"
Set-TextCell "D5" "281474991265569"
Set-TextCell "E5" "103"
Set-TextCell "F5" "52215661"
Set-TextCell "G5" "CARLOS ALBERTO JIMENEZ"
$ws.Range("H5").Value = 20.67005281
$ws.Range("I5").Value = -103.37155847
$ws.Range("J5").Value = 0.7675963044166565
Set-TextCell "K5" "No video URL"
Set-TextCell "L5" "No video URL"

# --- Row 6 -----------------------------------------------------------
Set-TextCell "A6" "281474991395157-1738864139858"
Set-TextCell "B6" "Harsh Brake"
Set-TextCell "C6" "2025-02-06T11:48:59.858"
Set-TextCell "D6" "281474991395157"
Set-TextCell "E6" "126"
Set-TextCell "F6" "No driver ID"
Set-TextCell "G6" "No driver name"
$ws.Range("H6").Value = 20.68131616
$ws.Range("I6").Value = -103.310700499
$ws.Range("J6").Value = 0.939757764339447
Set-TextCell "K6" "No video URL"
Set-TextCell "L6" "No video URL"

# --- Former row 7 is removed entirely ---------------------------------
$ws.Rows(7).Delete()
